# The sheet holds one row per year (2000..2020). The update drops the ten
# oldest years (2000..2009), keeps 2010..2020 as-is, and appends two more
# years: 2021 (complete) and 2022 (only "施工房屋面积" C and "竣工房屋面积" F
# are known so far - the other three measures are blank for that year).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop rows 2-11 (years 2000..2009) entirely so 2010..2020 (old rows 12-22)
# shift up to rows 2-12, and the sheet's used range shrinks accordingly.
$ws.Range("A2:F11").EntireRow.Delete()

# Rows 2-12 now hold 2010年..2020年 unchanged. Append the new years.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 10.4
$ws.Range("C13").Value = 975386.5101
$ws.Range("D13").Value = 39458.1541
$ws.Range("E13").Value = 3891
$ws.Range("F13").Value = 101411.9393

$ws.Range("A14").Value = "2022年"
# A bare quote-prefix (no text) stores as an empty, text-typed cell instead
# of clearing it outright, matching the "reported but value pending" blanks
# for B/D/E in the source data.
$ws.Range("B14").Value = "'"
$ws.Range("C14").Value = 904999.2574
$ws.Range("D14").Value = "'"
$ws.Range("E14").Value = "'"
$ws.Range("F14").Value = 86222.22100000001
# Drop the quote-prefix formatting flag so the blanks keep the sheet's
# regular (unstyled) number formatting, like every other non-year column.
$ws.Range("B14").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"

# Column A uses a bold/bordered/centered style; copy it onto the two new
# year labels so they match the rest of the column.
$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)
$excel.CutCopyMode = 0
